$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: round the Ost/Nord (Q/R) coordinates to whole numbers.
$ws.Range("Q3").Value = 446627
$ws.Range("R3").Value = 7032919

# Row 4 and Row 5 swap their species/record data (A,B,D,E,F,G,H), and the
# Ost/Nord coordinates for the swapped rows are rounded to whole numbers.
# (Use .Value2 -- .Value round-trips through a COM Variant descriptor
# instead of the scalar in this host, so read with .Value2.)
$row4A = $ws.Range("A4").Value2
$row4B = $ws.Range("B4").Value2
$row4D = $ws.Range("D4").Value2
$row4E = $ws.Range("E4").Value2
$row4F = $ws.Range("F4").Value2
$row4G = $ws.Range("G4").Value2
$row4H = $ws.Range("H4").Value2

$row5A = $ws.Range("A5").Value2
$row5B = $ws.Range("B5").Value2
$row5D = $ws.Range("D5").Value2
$row5E = $ws.Range("E5").Value2
$row5F = $ws.Range("F5").Value2
$row5G = $ws.Range("G5").Value2
$row5H = $ws.Range("H5").Value2

$ws.Range("A4").Value = $row5A
$ws.Range("B4").Value = $row5B
$ws.Range("D4").Value = $row5D
$ws.Range("E4").Value = $row5E
$ws.Range("F4").Value = $row5F
$ws.Range("G4").Value = $row5G
$ws.Range("H4").Value = $row5H
$ws.Range("Q4").Value = 446547
$ws.Range("R4").Value = 7032732

$ws.Range("A5").Value = $row4A
$ws.Range("B5").Value = $row4B
$ws.Range("D5").Value = $row4D
$ws.Range("E5").Value = $row4E
$ws.Range("F5").Value = $row4F
$ws.Range("G5").Value = $row4G
$ws.Range("H5").Value = $row4H
$ws.Range("Q5").Value = 446544
$ws.Range("R5").Value = 7032738

# Clear the Starttid/Sluttid (Z/AB) cells for rows 3-5; the diff removes
# them entirely rather than leaving blank string values.
$ws.Range("Z3:Z5").ClearContents()
$ws.Range("AB3:AB5").ClearContents()
